# Append: 2025-09-20 18:28 JST
# Update the "取得日時" (acquired datetime) column (A) for all existing
# data rows (2-12) on the active sheet ("ランサーズ") to the new
# timestamp recorded at the time of this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-09-20 18:28:43"

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
